# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 2 = "H")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 446
$wsOff.Range("C2").Value = 284
$wsOff.Range("D2").Value = 130
$wsOff.Range("E2").Value = 54

# Update DEF sheet (row 2 = "H")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 448
$wsDef.Range("C2").Value = 300
$wsDef.Range("D2").Value = 128
$wsDef.Range("E2").Value = 55
$wsDef.Range("F2").Value = 13
